{"js": "// Find the target run of text in the task-list paragraph and insert the\n// missing phrase \"\u0434\u0430\u043d\u0438\u0445 \u0434\u043b\u044f \u0432\u0438\u043a\u043e\u043d\u0430\u043d\u043d\u044f \u043e\u043f\u0435\u0440\u0430\u0446\u0456\u0439\" right after\n// \"\u0444\u0443\u043d\u043a\u0446\u0456\u0454\u044e \u0437\u0430\u043f\u0438\u0442\u0456\u0432\" and before \" \u0432\u0456\u0434 \u043a\u043e\u0440\u0438\u0441\u0442\u0443\u0432\u0430\u0447\u0430\".\nconst searchResults = context.document.body.search(\"\u0444\u0443\u043d\u043a\u0446\u0456\u0454\u044e \u0437\u0430\u043f\u0438\u0442\u0456\u0432 \u0432\u0456\u0434 \u043a\u043e\u0440\u0438\u0441\u0442\u0443\u0432\u0430\u0447\u0430\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target phrase not found in document.\");\n}\n\nconst hit = searchResults.items[0];\nhit.insertText(\n  \"\u0444\u0443\u043d\u043a\u0446\u0456\u0454\u044e \u0437\u0430\u043f\u0438\u0442\u0456\u0432 \u0434\u0430\u043d\u0438\u0445 \u0434\u043b\u044f \u0432\u0438\u043a\u043e\u043d\u0430\u043d\u043d\u044f \u043e\u043f\u0435\u0440\u0430\u0446\u0456\u0439 \u0432\u0456\u0434 \u043a\u043e\u0440\u0438\u0441\u0442\u0443\u0432\u0430\u0447\u0430\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Task 1 (topic 04): extend the sentence about the calculator's query\n# feature so it also mentions reading the data needed to perform the\n# operations, right before \"... from the user\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute(\n    \"\u0444\u0443\u043d\u043a\u0446\u0456\u0454\u044e \u0437\u0430\u043f\u0438\u0442\u0456\u0432 \u0432\u0456\u0434 \u043a\u043e\u0440\u0438\u0441\u0442\u0443\u0432\u0430\u0447\u0430\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"\u0444\u0443\u043d\u043a\u0446\u0456\u0454\u044e \u0437\u0430\u043f\u0438\u0442\u0456\u0432 \u0434\u0430\u043d\u0438\u0445 \u0434\u043b\u044f \u0432\u0438\u043a\u043e\u043d\u0430\u043d\u043d\u044f \u043e\u043f\u0435\u0440\u0430\u0446\u0456\u0439 \u0432\u0456\u0434 \u043a\u043e\u0440\u0438\u0441\u0442\u0443\u0432\u0430\u0447\u0430\",\n    2\n)\n\nif (-not $found) {\n    throw \"Target phrase not found in document.\"\n}\n"}
